# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Update "Hoja1" A1 text block with new conversion rates ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$oldText = $wsHoja1.Range("A1").Value()
$find = "✅ 1000 Bs = 4.39 = 17141.1 pesos`n✅ 17141.1 pesos = 4.36 = 942.78 Bs"
$replace = "✅ 1000 Bs = 4.49 = 17499.54 pesos`n✅ 17499.54 pesos = 4.45 = 948.48 Bs"
$newText = $oldText.Replace($find, $replace)
$wsHoja1.Range("A1").Value = $newText

# --- Update "tasas" sheet rate cells ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 222.895
$wsTasas.Range("O10").Value = 3900.56
$wsTasas.Range("N12").Value = 3929.85
$wsTasas.Range("O12").Value = 213
